$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-20 Saturday" "2024-01-21 Sunday"

Replace-Text "316×5=" "346×3="
Replace-Text "343×7=" "161×4="
Replace-Text "861×8=" "153×3="
Replace-Text "623×5=" "655×6="
Replace-Text "866×2=" "258×2="

Replace-Text "903×9=" "216×5="
Replace-Text "303×8=" "185×9="
Replace-Text "615×9=" "829×8="
Replace-Text "132×5=" "602×7="
Replace-Text "740×2=" "901×5="

Replace-Text "925×3=" "233×5="
Replace-Text "844×8=" "474×5="
Replace-Text "525×2=" "970×2="
Replace-Text "801×6=" "715×5="
Replace-Text "566×8=" "710×8="

Replace-Text "538×4=" "693×5="
Replace-Text "514×4=" "564×9="
Replace-Text "724×6=" "380×2="
Replace-Text "212×2=" "832×6="
Replace-Text "109×9=" "698×2="

Replace-Text "733×9=" "973×5="
Replace-Text "548×5=" "516×3="
Replace-Text "621×6=" "456×4="
Replace-Text "428×4=" "985×2="
Replace-Text "163×8=" "637×2="
